$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new student rows (6-9)
$ws.Range("A6").Value = "Danial"
$ws.Range("B6").Value = 18
$ws.Range("C6").Value = 77

$ws.Range("A7").Value = "Qaisara"
$ws.Range("B7").Value = 16
$ws.Range("C7").Value = 76

$ws.Range("A8").Value = "Raif"
$ws.Range("B8").Value = 14
$ws.Range("C8").Value = 99

$ws.Range("A9").Value = "Azfar"
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 75

# Match the final selection left after data entry
$ws.Range("A10").Select()
